# issue #5: add portion and total(area*portion) of land
# The source data rows (row 2) are pulled up into the header row (row 1),
# the old leading index column (A) is cleared, and the now-empty data row
# is removed on both the vehicle ("汽車") and debt ("債務") sheets.

$wb = $excel.ActiveWorkbook

foreach ($ws in @($wb.Worksheets.Item(1), $wb.Worksheets.Item(2))) {
    $b2 = $ws.Range("B2").Value()
    $c2 = $ws.Range("C2").Value()
    $d2 = $ws.Range("D2").Value()
    $e2 = $ws.Range("E2").Value()
    $f2 = $ws.Range("F2").Value()
    $g2 = $ws.Range("G2").Value()

    $ws.Range("B1").Value = $b2
    $ws.Range("C1").Value = $c2
    $ws.Range("D1").Value = $d2
    $ws.Range("E1").Value = $e2
    $ws.Range("F1").Value = $f2
    $ws.Range("G1").Value = $g2

    $ws.Rows.Item(2).Delete()
    $ws.Range("A1").Clear()
}
